# Refresh cryptos list snapshot (GitHub Actions scheduled update).
# Prices/volumes are stored as plain text (not numbers), so we force the
# "Text" number format before assigning, then restore the default "Normal"
# style so no stray formatting is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "41.675.10"
Set-TextValue $ws.Range("E2") "  +5.28%  "
Set-TextValue $ws.Range("D3") "2.223.26"
Set-TextValue $ws.Range("E3") "  +3.23%  "
Set-TextValue $ws.Range("E4") "  +0.04%  "
Set-TextValue $ws.Range("D5") "230.86"
Set-TextValue $ws.Range("E5") "  +1.84%  "
Set-TextValue $ws.Range("E6") "  +0.80%  "
Set-TextValue $ws.Range("D7") "61.00"
Set-TextValue $ws.Range("E7") "  -2.51%  "
Set-TextValue $ws.Range("E8") "  +0.06%  "
Set-TextValue $ws.Range("D9") "0.401"
Set-TextValue $ws.Range("E9") "  +3.27%  "
Set-TextValue $ws.Range("D10") "58.75"
Set-TextValue $ws.Range("E10") "  +0.65%  "
Set-TextValue $ws.Range("D11") "0.0890"
Set-TextValue $ws.Range("E11") "  +5.93%  "
Set-TextValue $ws.Range("E12") "  +0.30%  "
Set-TextValue $ws.Range("D13") "2.553.77"
Set-TextValue $ws.Range("E13") "  +3.21%  "
Set-TextValue $ws.Range("D14") "15.63"
Set-TextValue $ws.Range("E14") "  -0.97%  "
Set-TextValue $ws.Range("E15") "  +0.17%  "
Set-TextValue $ws.Range("D16") "0.797"
Set-TextValue $ws.Range("E16") "  -0.65%  "
Set-TextValue $ws.Range("E17") "  +1.62%  "
Set-TextValue $ws.Range("D18") "2.235.47"
Set-TextValue $ws.Range("E18") "  +3.98%  "
Set-TextValue $ws.Range("D19") "41.505.17"
Set-TextValue $ws.Range("E19") "  +5.03%  "
Set-TextValue $ws.Range("D20") "72.80"
Set-TextValue $ws.Range("E20") "  +1.71%  "
Set-TextValue $ws.Range("D21") "0.0₃0897"
Set-TextValue $ws.Range("E21") "  +5.78%  "
Set-TextValue $ws.Range("D22") "6.04"
Set-TextValue $ws.Range("E22") "  -0.42%  "
Set-TextValue $ws.Range("D23") "249.84"
Set-TextValue $ws.Range("E23") "  +10.03%  "
Set-TextValue $ws.Range("E24") "  +0.01%  "
Set-TextValue $ws.Range("D25") "2.40"
Set-TextValue $ws.Range("E25") "  +2.30%  "
Set-TextValue $ws.Range("E26") "  +2.78%  "
Set-TextValue $ws.Range("E27") "  +1.88%  "
Set-TextValue $ws.Range("D28") "167.90"
Set-TextValue $ws.Range("E28") "  -1.61%  "
Set-TextValue $ws.Range("E29") "  +1.55%  "
Set-TextValue $ws.Range("E30") "  +2.02%  "
Set-TextValue $ws.Range("E31") "  -0.88%  "
Set-TextValue $ws.Range("E32") "  -2.62%  "
Set-TextValue $ws.Range("E33") "  +0.79%  "
Set-TextValue $ws.Range("D34") "4.94"
Set-TextValue $ws.Range("E34") "  +5.14%  "
Set-TextValue $ws.Range("E35") "  +1.20%  "
Set-TextValue $ws.Range("E36") "  +1.43%  "
Set-TextValue $ws.Range("D37") "6.56"
Set-TextValue $ws.Range("E37") "  -5.55%  "
Set-TextValue $ws.Range("D38") "3.67"
Set-TextValue $ws.Range("E38") "  -2.21%  "
Set-TextValue $ws.Range("D39") "2.35"
Set-TextValue $ws.Range("E39") "  -1.41%  "
Set-TextValue $ws.Range("D40") "0.000245"
Set-TextValue $ws.Range("E40") "  +28.44%  "
Set-TextValue $ws.Range("E41") "  -0.09%  "
Set-TextValue $ws.Range("D42") "4.87"
Set-TextValue $ws.Range("E42") "  +1.51%  "
Set-TextValue $ws.Range("D43") "0.0238"
Set-TextValue $ws.Range("E43") "  +5.09%  "
Set-TextValue $ws.Range("D44") "8.58"
Set-TextValue $ws.Range("E44") "  +9.22%  "
Set-TextValue $ws.Range("D45") "0.0980"
Set-TextValue $ws.Range("E45") "  +6.74%  "
Set-TextValue $ws.Range("E46") "  +2.00%  "
Set-TextValue $ws.Range("D47") "98.94"
Set-TextValue $ws.Range("E47") "  -3.92%  "
Set-TextValue $ws.Range("D48") "1.463.56"
Set-TextValue $ws.Range("E48") "  -3.23%  "
Set-TextValue $ws.Range("B49") "HuobiToken"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D49") "2.80"
Set-TextValue $ws.Range("E49") "  -0.04%  "
Set-TextValue $ws.Range("B50") "InjectiveProtocol"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D50") "16.42"
Set-TextValue $ws.Range("E50") "  -6.68%  "
